# Update view-count figures (column F) on several sheets to reflect the
# latest generated output (gh-pages rebuild at commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 21230
$ws1.Range("F5").Value  = 3036
$ws1.Range("F6").Value  = 809
$ws1.Range("F14").Value = 524
$ws1.Range("F16").Value = 270

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6135
$ws3.Range("F5").Value = 1636

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 6135
$ws4.Range("F5").Value  = 1636
$ws4.Range("F7").Value  = 21230
$ws4.Range("F13").Value = 3037
$ws4.Range("F14").Value = 809
$ws4.Range("F29").Value = 524
$ws4.Range("F33").Value = 270
